$d = $word.ActiveDocument

# Step 1: Text replace to get the ECOG expansion
$d.Content.Find.Execute(
    "She assesses his performance status as ECOG 1. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "She assesses his performance status as Eastern Cooperative Oncology Group (ECOG) 1. ",
    2) | Out-Null

Write-Output $d.Content.Text
